$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "adelina41736"
$ws.Range("B1").Value = "Maricruz"
$ws.Range("C1").Value = "Zamora"
$ws.Range("D1").Value = "igalindo@gmail.com"
$ws.Range("E1").Value = "W5V&vdsE0LHZ"
$ws.Range("F1").Value = "W5V&vdsE0LHZ"
$ws.Range("G1").Value = "Válido"

# Row 2
$ws.Range("A2").Value = "samanta15396"
$ws.Range("B2").Value = "Ruth"
$ws.Range("C2").Value = "Silva"
$ws.Range("D2").Value = "flastra@inversiones.es"
$ws.Range("E2").Value = "m7FUa5ed*CIF"
$ws.Range("F2").Value = "m7FUa5ed*CIF"
$ws.Range("G2").Value = "Válido"

# Row 3 (new)
$ws.Range("A3").Value = "maciassarita605"
$ws.Range("B3").Value = "Elena"
$ws.Range("C3").Value = "Solana"
$ws.Range("D3").Value = "emilia67@gmail.com"
$ws.Range("E3").Value = "2Tq(uBue#^o1"
$ws.Range("F3").Value = "2Tq(uBue#^o1"
$ws.Range("G3").Value = "Válido"
